$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.377.11'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -1.83%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.650.59'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -3.38%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.92'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3646'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -2.87%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '46.95'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -5.43%  '
$ws.Range("E9").Value = '  -5.81%  '
$ws.Range("E10").Value = '  -7.26%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07019'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -6.99%  '
$ws.Range("E12").Value = '  +0.13%  '
$ws.Range("E13").Value = '  -5.93%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '19.33'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -8.36%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.583'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -6.49%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.653.27'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -3.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001040'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -8.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06593'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.03%  '
$ws.Range("E19").Value = '  +0.18%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '77.92'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -7.74%  '
$ws.Range("E21").Value = '  -7.26%  '
$ws.Range("E22").Value = '  -10.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.41'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -6.25%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.373.26'
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.480'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.330'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -16.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '147.05'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -3.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.53'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -9.11%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.834.33'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -3.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.66'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -6.74%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.164'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -6.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.069'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.65%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.624'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -18.97%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08444'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.83%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.668'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -8.90%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '12.24'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -11.18%  '
$ws.Range("E37").Value = '  -7.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.245'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06017'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -9.84%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.02209'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -8.28%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.2060'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -7.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.116'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -13.04%  '
$ws.Range("E43").Value = '  +0.21%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5866'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -8.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.773'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.33%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '12.60'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -9.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5599'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -9.20%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '122.27'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -6.05%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.939'
$ws.Range("D49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06885'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -5.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '74.48'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -6.65%  '
